$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-1.239354553595497, 0.5934752043385932, 1.315698351103824),
    @(0.1265866749218411, 0.9830147470409047, 1.104435794153933),
    @(0.2933733014283093, 0.4318224985956808, 1.679565984772173),
    @(-0.3514636487162646, 1.303912767754777, 1.717408596215416),
    @(-0.8439442796328934, 1.107491066802433, -0.297189421712284),
    @(1.054603315660084, -0.3611965101983832, -2.38295376773451),
    @(1.400082142342559, -1.373465267046946, -4.113783430948127),
    @(2.811405012786107, 2.239890750284216, -2.48325762769725),
    @(1.144158274066543, 0.4079549139291005, -2.216124885407739),
    @(-3.973825000981487, -4.923421618171759, 5.359984481912402),
    @(-10.63765511113626, -7.195822839694799, 9.095965221589761),
    @(3.433085498305679, 4.357403949493808, -2.281185732013872),
    @(4.477401155732298, -3.086138354523986, -4.067132330150736),
    @(5.366742232297456, 3.331075757610629, 1.803214041672062),
    @(0.746474412569362, 1.614592228954578, -2.489844009214525),
    @(-4.899446143452909, 0.4346358306082141, 4.76701935377427),
    @(-9.7976465981438, 3.087984385994603, -1.579483953341611),
    @(-0.5613160348673585, 0.9123670508682955, 0.6450322871691367),
    @(2.186993496270989, -7.507715720437281, 1.97455261667396),
    @(2.968808529135423, 8.113890336998949, -6.538449691780965),
    @(-1.87808840600416, -1.777517291417992, 1.215824003261815),
    @(-4.240441676278492, 2.87217904187515, 8.334868672660862),
    @(1.2107466585311, 12.91206166397641, -4.052727646764785),
    @(-0.8707238068140057, -7.102825836988711, 0.5642036477901016),
    @(6.469962025529018, -3.154151007992533, -3.828252468865385),
    @(-1.114345516927603, -1.64307104963562, 0.5952699299951911),
    @(-2.965141663992466, -1.430339549081439, 1.28920519351959),
    @(-4.759042248326804, -5.845142855518262, -5.263579006762336),
    @(3.102689377537001, -4.462434116439107, 5.959941364070088),
    @(3.054766318871562, -4.391178798045361, 1.19998335943338)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
